$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 8-15 (B/C/D/E values), B text shifts automatically
# once the underlying strings are set directly (labels renamed to reflect
# 2 new "line" entries pushed in ahead of the "extr" series).

# Row 8 (id 6): extr1 -> line7
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9 (id 7): extr2 -> line8
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Row 10 (id 8): extr3 -> extr1
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11 (id 9): extr4 -> extr2
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12 (id 10): extr5 -> extr3
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

# Row 13 (id 11): extr6 -> extr4
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# Row 14 (id 12): extr7 -> extr5
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

# Row 15 (id 13): extr8 -> extr6
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# New row 16 (id 14): extr7
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(16, 1).PasteSpecial(-4122)
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# New row 17 (id 15): extr8
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial(-4122)
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $true

$excel.CutCopyMode = $false
